# Add new columns I ("I0") and J ("IF") to Sheet1, mirroring the
# existing header style (copied from H1) and filling in the data rows
# (2-79) with the values supplied in the commit's diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1): copy H1's formatting onto I1/J1, then set text ---
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("J1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-79: column I and J values ---
$iVals = @(3,9,8,7,1,9,1,1,7,6,6,8,8,6,9,1,1,7,1,8,6,7,8,8,8,6,8,8,5,8,8,8,8,8,4,5,6,6,8,8,10,11,7,1,9,9,9,6,9,6,8,7,7,9,7,9,7,8,6,6,8,8,8,8,8,6,4,7,7,9,8,6,2,5,7,5,4,4)
$jVals = @(5,9,8,7,1,9,1,1,7,7,6,8,8,6,9,1,1,7,1,8,6,7,8,8,8,6,8,8,6,8,8,8,8,8,4,5,7,6,8,8,10,11,7,1,9,9,9,6,9,6,8,7,7,9,8,9,7,9,6,6,8,8,8,8,8,6,5,7,7,9,8,6,2,5,7,5,4,4)

for ($i = 0; $i -lt $iVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$i]
    $ws.Cells.Item($row, 10).Value = $jVals[$i]
}
